{"js": "// Update the date line and all 25 division problems/answers in the table\n// to the new values described by the commit's regenerated output.\nconst replacements = [\n  [\"2025-08-17 Sunday\", \"2025-08-18 Monday\"],\n  [\"934\u00f79=103, 7\", \"792\u00f78=99, 0\"],\n  [\"491\u00f72=245, 1\", \"933\u00f72=466, 1\"],\n  [\"257\u00f72=128, 1\", \"711\u00f77=101, 4\"],\n  [\"380\u00f76=63, 2\", \"858\u00f75=171, 3\"],\n  [\"751\u00f79=83, 4\", \"737\u00f76=122, 5\"],\n  [\"623\u00f78=77, 7\", \"640\u00f72=320, 0\"],\n  [\"515\u00f72=257, 1\", \"727\u00f78=90, 7\"],\n  [\"581\u00f74=145, 1\", \"176\u00f75=35, 1\"],\n  [\"299\u00f76=49, 5\", \"767\u00f77=109, 4\"],\n  [\"986\u00f77=140, 6\", \"171\u00f73=57, 0\"],\n  [\"941\u00f74=235, 1\", \"693\u00f79=77, 0\"],\n  [\"948\u00f79=105, 3\", \"149\u00f72=74, 1\"],\n  [\"122\u00f77=17, 3\", \"630\u00f79=70, 0\"],\n  [\"490\u00f73=163, 1\", \"709\u00f75=141, 4\"],\n  [\"816\u00f74=204, 0\", \"866\u00f79=96, 2\"],\n  [\"879\u00f73=293, 0\", \"660\u00f72=330, 0\"],\n  [\"249\u00f74=62, 1\", \"557\u00f78=69, 5\"],\n  [\"897\u00f77=128, 1\", \"810\u00f78=101, 2\"],\n  [\"477\u00f72=238, 1\", \"914\u00f76=152, 2\"],\n  [\"683\u00f77=97, 4\", \"191\u00f76=31, 5\"],\n  [\"354\u00f73=118, 0\", \"987\u00f77=141, 0\"],\n  [\"224\u00f75=44, 4\", \"438\u00f74=109, 2\"],\n  [\"721\u00f72=360, 1\", \"829\u00f74=207, 1\"],\n  [\"717\u00f79=79, 6\", \"993\u00f73=331, 0\"],\n  [\"541\u00f75=108, 1\", \"877\u00f74=219, 1\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Not found: \" + oldText);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and all 25 division problems/answers in the table\n# to the new values described by the commit's regenerated output.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2025-08-17 Sunday\", \"2025-08-18 Monday\"),\n  @(\"934\u00f79=103, 7\", \"792\u00f78=99, 0\"),\n  @(\"491\u00f72=245, 1\", \"933\u00f72=466, 1\"),\n  @(\"257\u00f72=128, 1\", \"711\u00f77=101, 4\"),\n  @(\"380\u00f76=63, 2\", \"858\u00f75=171, 3\"),\n  @(\"751\u00f79=83, 4\", \"737\u00f76=122, 5\"),\n  @(\"623\u00f78=77, 7\", \"640\u00f72=320, 0\"),\n  @(\"515\u00f72=257, 1\", \"727\u00f78=90, 7\"),\n  @(\"581\u00f74=145, 1\", \"176\u00f75=35, 1\"),\n  @(\"299\u00f76=49, 5\", \"767\u00f77=109, 4\"),\n  @(\"986\u00f77=140, 6\", \"171\u00f73=57, 0\"),\n  @(\"941\u00f74=235, 1\", \"693\u00f79=77, 0\"),\n  @(\"948\u00f79=105, 3\", \"149\u00f72=74, 1\"),\n  @(\"122\u00f77=17, 3\", \"630\u00f79=70, 0\"),\n  @(\"490\u00f73=163, 1\", \"709\u00f75=141, 4\"),\n  @(\"816\u00f74=204, 0\", \"866\u00f79=96, 2\"),\n  @(\"879\u00f73=293, 0\", \"660\u00f72=330, 0\"),\n  @(\"249\u00f74=62, 1\", \"557\u00f78=69, 5\"),\n  @(\"897\u00f77=128, 1\", \"810\u00f78=101, 2\"),\n  @(\"477\u00f72=238, 1\", \"914\u00f76=152, 2\"),\n  @(\"683\u00f77=97, 4\", \"191\u00f76=31, 5\"),\n  @(\"354\u00f73=118, 0\", \"987\u00f77=141, 0\"),\n  @(\"224\u00f75=44, 4\", \"438\u00f74=109, 2\"),\n  @(\"721\u00f72=360, 1\", \"829\u00f74=207, 1\"),\n  @(\"717\u00f79=79, 6\", \"993\u00f73=331, 0\"),\n  @(\"541\u00f75=108, 1\", \"877\u00f74=219, 1\")\n)\n\nforeach ($pair in $pairs) {\n  $old = $pair[0]\n  $new = $pair[1]\n\n  $r = $d.Content\n  $r.Find.ClearFormatting()\n  $r.Find.Replacement.ClearFormatting()\n  $r.Find.Text = $old\n  $r.Find.Replacement.Text = $new\n  $found = $r.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n  if (-not $found) {\n    throw \"Not found: $old\"\n  }\n}\n"}
